# Natmi following Dr Hou advice
# Update LR-pair rows for Egf-Egfr (OldD0) with the recomputed specificity values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,("ECs", "Egf", "Egfr", "ECs", 1, 0.3333333333333333, 0.050561, 0.151683, 0.1845256053410153, 0.1845256053410153, 3, 1, 2.08532, 6.25596, 0.01753772176136817, 0.01753772176136816, 0.10543586452, 0.94892278068, 0.003236158724318758, 0.003236158724318757)
    ,("ECs", "Egf", "Egfr", "FAPs", 1, 0.3333333333333333, 0.050561, 0.151683, 0.1845256053410153, 0.1845256053410153, 3, 1, 101.898173, 305.694519, 0.8569724579756384, 0.8569724579756383, 5.152073525053, 46.368661725477, 0.1581333615685324, 0.1581333615685324)
    ,("ECs", "Egf", "Egfr", "sCs", 1, 0.3333333333333333, 0.050561, 0.151683, 0.1845256053410153, 0.1845256053410153, 3, 1, 14.921347, 44.76404100000001, 0.1254898202629935, 0.1254898202629935, 0.7544382256670001, 6.789944031003001, 0.02315608504816408, 0.02315608504816408)
    ,("FAPs", "Egf", "Egfr", "ECs", 3, 1, 0.1101723333333333, 0.3305169999999999, 0.4020809813920896, 0.4020809813920896, 3, 1, 2.08532, 6.25596, 0.01753772176136817, 0.01753772176136816, 0.2297445701466666, 2.06770113132, 0.007051584377192319, 0.007051584377192318)
    ,("FAPs", "Egf", "Egfr", "FAPs", 3, 1, 0.1101723333333333, 0.3305169999999999, 0.4020809813920896, 0.4020809813920896, 3, 1, 101.898173, 305.694519, 0.8569724579756384, 0.8569724579756383, 11.22635948181366, 101.037235336323, 0.344572326928836, 0.3445723269288359)
    ,("FAPs", "Egf", "Egfr", "sCs", 3, 1, 0.1101723333333333, 0.3305169999999999, 0.4020809813920896, 0.4020809813920896, 3, 1, 14.921347, 44.76404100000001, 0.1254898202629935, 0.1254898202629935, 1.643919615466333, 14.795276539197, 0.05045707008606137, 0.05045707008606136)
    ,("sCs", "Egf", "Egfr", "ECs", 3, 1, 0.113272, 0.339816, 0.413393413266895, 0.413393413266895, 3, 1, 2.08532, 6.25596, 0.01753772176136817, 0.01753772176136816, 0.23620836704, 2.12587530336, 0.007249978659857089, 0.007249978659857087)
    ,("sCs", "Egf", "Egfr", "FAPs", 3, 1, 0.113272, 0.339816, 0.413393413266895, 0.413393413266895, 3, 1, 101.898173, 305.694519, 0.8569724579756384, 0.8569724579756383, 11.542209852056, 103.879888668504, 0.3542667694782699, 0.3542667694782698)
    ,("sCs", "Egf", "Egfr", "sCs", 3, 1, 0.113272, 0.339816, 0.413393413266895, 0.413393413266895, 3, 1, 14.921347, 44.76404100000001, 0.1254898202629935, 0.1254898202629935, 1.690170817384, 15.211537356456, 0.05187666512876806, 0.05187666512876805)
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range($cols[$j] + $rowNum).Value2 = $rowVals[$j]
    }
}
